$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the rows that changed
$ws.Range("F2").Value = -5
$ws.Range("F4").Value = -4
$ws.Range("F6").Value = 5
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = -3
$ws.Range("F13").Value = -8
$ws.Range("F16").Value = 10
$ws.Range("F18").Value = -2
$ws.Range("F21").Value = -2
$ws.Range("F28").Value = 3
